# Update forecast figures for VentaNueva Sheet1 and append the new day (43521).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Revised figures for existing rows (International = B, National = C) ---
$ws.Range("B18").Value = 514403
$ws.Range("C18").Value = 2605004

$ws.Range("B19").Value = 648503

$ws.Range("B20").Value = 1359406

$ws.Range("B22").Value = 668072
$ws.Range("C22").Value = 3519197

$ws.Range("B23").Value = 392336

$ws.Range("C24").Value = 5964662

$ws.Range("B25").Value = 668969
$ws.Range("C25").Value = 2381570

# --- Grow the table to include the new row 26, then fill it in ---
$lo = $ws.ListObjects.Item("Tabla1")
$lo.Resize($ws.Range("A1:D26"))

$ws.Range("A25").Copy()
$ws.Range("A26").PasteSpecial(-4122)

$ws.Range("A26").Value = 43521
$ws.Range("B26").Value = 1422652
$ws.Range("C26").Value = 4428691
$ws.Range("D26").Formula = "=B26+C26/Hoja2!`$A`$2"
